# Swap the "codeforiati:group-name" (column D) and "codeforiati:group-code"
# (column E) columns - both their header labels and every data value - on
# the ReportingOrganisationGroup sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row based on column A.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 1; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2

    $dCell.Value2 = $eVal
    $eCell.Value2 = $dVal
}
